# Edit: insert 4 new price-report rows (Mandarina - Clemenuless) into the
# "Macroferia Regional de Talca" daily price sheet, right before the
# existing row 51, shifting the remaining rows (old 51-116) down to 55-120.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at position 51; this pushes old rows 51-116 down to 55-120
$ws.Rows("51:54").Insert()

# NOTE: this engine's PowerShell does not bind named (-Param value) arguments
# to function parameters, so positional parameters are used instead.
function Set-PrecioRow {
    param($Row, $Fecha, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $Origen, $PrecioKg, $KgUnidad)

    $ws.Cells.Item($Row, 1).Value  = 5
    $ws.Cells.Item($Row, 2).Value  = "Macroferia Regional de Talca"
    $ws.Cells.Item($Row, 3).Value  = "Maule"
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 5).Value  = 7
    $ws.Cells.Item($Row, 6).Value  = "Fruta"
    $ws.Cells.Item($Row, 7).Value  = 100102
    $ws.Cells.Item($Row, 8).Value  = "Cítricos"
    $ws.Cells.Item($Row, 9).Value  = 100102004
    $ws.Cells.Item($Row, 10).Value = "Mandarina"
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $KgUnidad
}

# New row 51: Clemenuless / Primera
Set-PrecioRow 51 44413 "Clemenuless" "Primera" 210 6000 6000 6000 "$/bandeja 10 kilos" "Provincia de Limarí" 600 10

# New row 52: Clemenuless / Primera
Set-PrecioRow 52 44413 "Clemenuless" "Primera" 150 8000 8000 8000 "$/caja 18 kilos" "Provincia de Limarí" 444 18

# New row 53: Clemenuless / Segunda
Set-PrecioRow 53 44413 "Clemenuless" "Segunda" 110 6000 6000 6000 "$/caja 18 kilos" "Provincia de Limarí" 333 18

# New row 54: Clemenuless / Tercera
Set-PrecioRow 54 44413 "Clemenuless" "Tercera" 100 5000 5000 5000 "$/caja 18 kilos" "Provincia de Limarí" 278 18
